$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3046.375
$ws.Range("J32").Value = 3095.5
$ws.Range("L32").Value = 3095.5
$ws.Range("N32").Value = -3747.5

$ws.Range("H33").Value = 252.6923
$ws.Range("I33").Value = 171.36363
$ws.Range("K33").Value = 171.36363
$ws.Range("M33").Value = 57.63637

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = ""

$ws.Range("H98").Value = 998
$ws.Range("I98").Value = 998
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 998
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 500
$ws.Range("N98").Value = ""

$ws.Range("H116").Value = 3751.9048
$ws.Range("I116").Value = 3089
$ws.Range("J116").Value = 4354.5454
$ws.Range("K116").Value = 3089
$ws.Range("L116").Value = 4354.5454
$ws.Range("M116").Value = 353
$ws.Range("N116").Value = -11238.5454

$ws.Range("H122").Value = 998
$ws.Range("I122").Value = 998
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2994
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -544
$ws.Range("N122").Value = ""

$ws.Range("H129").Value = 1947.091
$ws.Range("I129").Value = 2245
$ws.Range("J129").Value = 1698.8334
$ws.Range("K129").Value = 6735
$ws.Range("L129").Value = 5096.5002
$ws.Range("M129").Value = -1735
$ws.Range("N129").Value = -15096.5002

$ws.Range("H132").Value = 1827.56
$ws.Range("I132").Value = 1732.4546
$ws.Range("K132").Value = 5197.3638
$ws.Range("M132").Value = -2667.3638

$ws.Range("H135").Value = 1285.7858
$ws.Range("I135").Value = 1285.7858
$ws.Range("K135").Value = 11572.0722
$ws.Range("M135").Value = -9037.072200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 420
$ws.Range("I22").Value = 425
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 425
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -126
$ws.Range("N22").Value = -998

$ws.Range("H26").Value = 32500
$ws.Range("I26").Value = 32500
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 32500
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -32170
$ws.Range("N26").Value = ""

$ws.Range("H45").Value = 1847.0555
$ws.Range("I45").Value = 1233.5
$ws.Range("J45").Value = 3994.5
$ws.Range("K45").Value = 1233.5
$ws.Range("L45").Value = 3994.5
$ws.Range("M45").Value = -856.5
$ws.Range("N45").Value = -4748.5

$ws.Range("H46").Value = 9524.5
$ws.Range("I46").Value = 9574
$ws.Range("J46").Value = 9475
$ws.Range("K46").Value = 9574
$ws.Range("L46").Value = 9475
$ws.Range("M46").Value = -9255
$ws.Range("N46").Value = -10113

$ws.Range("H63").Value = 6169.1
$ws.Range("I63").Value = 4770.143
$ws.Range("K63").Value = 4770.143
$ws.Range("M63").Value = -4084.143

$ws.Range("H66").Value = 6169.1
$ws.Range("I66").Value = 4770.143
$ws.Range("K66").Value = 23850.715
$ws.Range("M66").Value = -20418.715

$ws.Range("H132").Value = 1470.8334
$ws.Range("I132").Value = 1277.2858
$ws.Range("K132").Value = 3831.8574
$ws.Range("M132").Value = -1301.8574

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 16496
$ws.Range("I33").Value = 1000
$ws.Range("J33").Value = 21661.334
$ws.Range("K33").Value = 1000
$ws.Range("L33").Value = 21661.334
$ws.Range("M33").Value = -664
$ws.Range("N33").Value = -22333.334

$ws.Range("H86").Value = 5751.5
$ws.Range("I86").Value = 5751.5
$ws.Range("K86").Value = 5751.5
$ws.Range("M86").Value = -4628.5

$ws.Range("H89").Value = 5751.5
$ws.Range("I89").Value = 5751.5
$ws.Range("K89").Value = 28757.5
$ws.Range("M89").Value = -23141.5

$ws.Range("H94").Value = 1581.1666
$ws.Range("I94").Value = 1557.3334
$ws.Range("J94").Value = 1748
$ws.Range("K94").Value = 1557.3334
$ws.Range("L94").Value = 1748
$ws.Range("M94").Value = -1106.3334
$ws.Range("N94").Value = -2650

$ws.Range("H107").Value = 4377.909
$ws.Range("I107").Value = 4057.25
$ws.Range("K107").Value = 4057.25
$ws.Range("M107").Value = -2137.25

$ws.Range("H134").Value = 3763.5833
$ws.Range("I134").Value = 3763.5833
$ws.Range("K134").Value = 11290.7499
$ws.Range("M134").Value = -8755.749899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 49999
$ws.Range("J106").Value = 49999
$ws.Range("L106").Value = 49999
$ws.Range("N106").Value = -52523

$ws.Range("H122").Value = 1030.4828
$ws.Range("I122").Value = 1052.7059
$ws.Range("K122").Value = 3158.1177
$ws.Range("M122").Value = -708.1176999999998

$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800

$ws.Range("H132").Value = 4554.316
$ws.Range("I132").Value = 4554.316
$ws.Range("K132").Value = 13662.948
$ws.Range("M132").Value = -11132.948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 507.75
$ws.Range("I13").Value = 10.333333
$ws.Range("K13").Value = 30.999999
$ws.Range("M13").Value = 137.000001

$ws.Range("H68").Value = 4117.069
$ws.Range("J68").Value = 4125.5356
$ws.Range("L68").Value = 12376.6068
$ws.Range("N68").Value = -13998.6068

$ws.Range("H71").Value = 4117.069
$ws.Range("J71").Value = 4125.5356
$ws.Range("L71").Value = 37129.8204
$ws.Range("N71").Value = -45241.8204

$ws.Range("H107").Value = 1385.2727
$ws.Range("J107").Value = 1378.4
$ws.Range("L107").Value = 4135.200000000001
$ws.Range("N107").Value = -7975.200000000001

$ws.Range("H137").Value = 3247.25
$ws.Range("J137").Value = 3531.3333
$ws.Range("L137").Value = 10593.9999
$ws.Range("N137").Value = -20793.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0

$ws.Range("H122").Value = 2811.7144
$ws.Range("J122").Value = 2449.5
$ws.Range("L122").Value = 7348.5
$ws.Range("N122").Value = -12248.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7097.7
$ws.Range("I7").Value = 5796
$ws.Range("K7").Value = 5796
$ws.Range("M7").Value = -5684

$ws.Range("H22").Value = 869.8889
$ws.Range("I22").Value = 835.5714
$ws.Range("J22").Value = 990
$ws.Range("K22").Value = 835.5714
$ws.Range("L22").Value = 990
$ws.Range("M22").Value = -540.5714
$ws.Range("N22").Value = -1580

$ws.Range("H27").Value = 869.8889
$ws.Range("I27").Value = 835.5714
$ws.Range("J27").Value = 990
$ws.Range("K27").Value = 835.5714
$ws.Range("L27").Value = 990
$ws.Range("M27").Value = -728.5714
$ws.Range("N27").Value = -1204

$ws.Range("H46").Value = 1656.95
$ws.Range("J46").Value = 3577.75
$ws.Range("L46").Value = 3577.75
$ws.Range("N46").Value = -3953.75

$ws.Range("H82").Value = 681.1667
$ws.Range("I82").Value = 718.6
$ws.Range("K82").Value = 718.6
$ws.Range("M82").Value = -357.6

$ws.Range("H85").Value = 681.1667
$ws.Range("I85").Value = 718.6
$ws.Range("K85").Value = 718.6
$ws.Range("M85").Value = 529.4

$ws.Range("H100").Value = 5270.857
$ws.Range("I100").Value = 5270.857
$ws.Range("K100").Value = 5270.857
$ws.Range("M100").Value = -4729.857

$ws.Range("H104").Value = 11500
$ws.Range("J104").Value = 11500
$ws.Range("L104").Value = 11500
$ws.Range("N104").Value = -18488

$ws.Range("H126").Value = 7097.7
$ws.Range("I126").Value = 5796
$ws.Range("K126").Value = 17388
$ws.Range("M126").Value = -14918

$ws.Range("H135").Value = 94000
$ws.Range("J135").Value = 94000
$ws.Range("L135").Value = 94000
$ws.Range("N135").Value = -104140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17363.25
$ws.Range("I41").Value = 14989
$ws.Range("K41").Value = 14989
$ws.Range("M41").Value = -14599

$ws.Range("H113").Value = 420.8
$ws.Range("I113").Value = 387.57144
$ws.Range("J113").Value = 498.33334
$ws.Range("K113").Value = 1162.71432
$ws.Range("L113").Value = 1495.00002
$ws.Range("M113").Value = 1007.28568
$ws.Range("N113").Value = -5835.000019999999

$ws.Range("H122").Value = 3198.3667
$ws.Range("I122").Value = 2496.8948
$ws.Range("J122").Value = 4410
$ws.Range("K122").Value = 7490.6844
$ws.Range("L122").Value = 13230
$ws.Range("M122").Value = -5040.6844
$ws.Range("N122").Value = -18130

$ws.Range("H132").Value = 1811.5
$ws.Range("I132").Value = 1927.4286
$ws.Range("K132").Value = 5782.2858
$ws.Range("M132").Value = -3252.2858

$ws.Range("H136").Value = 2606.92
$ws.Range("I136").Value = 2606.92
$ws.Range("K136").Value = 7820.76
$ws.Range("M136").Value = -5270.76
